# ---------------------------------------------------------------------------
# Survey/Response ("SurveyRespnod") sheet update:
#  * Survey timeout raised to 8s -> big-message ipc/tcp survey runs that used
#    to time out now complete, so a whole new "ipc:" transport block (cols J-Q)
#    of big/small surveys-per-sec + KBs is captured, plus a new 8192 size row
#    and bigger sizes up to 1048576 that start failing once limits are hit.
#  * A second full run (rows 29-41) was captured underneath the first.
# ---------------------------------------------------------------------------

function Row2D {
    param([object[]]$vals)
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) { $arr[0,$i] = $vals[$i] }
    $arr
}
function Col2D {
    param([object[]]$vals)
    $n = $vals.Count
    $arr = New-Object 'object[,]' $n,1
    for ($i = 0; $i -lt $n; $i++) { $arr[$i,0] = $vals[$i] }
    $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SurveyRespnod")

# Row 3 headers: re-align N:Q block to big/big/small/small pattern, drop R3, add S3
$ws.Range("N3:Q3").Value2 = (Row2D @("big surveys/sec", "big KBs", "small surveys/sec", "small KBs"))
$ws.Range("R3").ClearContents()
$ws.Range("S3").Value2 = " "

# --- First run (rows 4-14) ---
$ws.Range("A4:A11").Value2 = (Col2D @(1024, 2048, 4096, 8192, 16384, 32768, 65536, 131072))

$ws.Range("B4:Q4").Value2 = (Row2D @(2614, 2164, 2472, 2472, 1422, 1422, 1432, 2863, 992, 992, 809, 2427, 682, 682, 749, 2998))
$ws.Range("B5:Q5").Value2 = (Row2D @(2419, 4838, 2082, 4163, 1380, 2759, 1345, 5380, 1001, 2001, 910, 5462, 700, 1399, 713, 5701))
$ws.Range("B6:Q6").Value2 = (Row2D @(2458, 9830, 1911, 7644, 1367, 5467, 1200, 9602, 943, 3771, 878, 10530, 751, 3003, 672, 10760))
$ws.Range("B7:Q7").Value2 = (Row2D @(2602, 20817, 1911, 15288, 1263, 10101, 1038, 16613, 974, 7796, 899, 21581, 681, 5447, 724, 23187))
$ws.Range("B8:Q8").Value2 = (Row2D @(2064, 33017, 1596, 25543, 1371, 21939, 1332, 42633, 955, 15286, 822, 39483, 697, 11151, 702, 44909))
$ws.Range("B9:Q9").Value2 = (Row2D @(2209, 70703, 1326, 42418, 1357, 43437, 997, 63783, 957, 30616, 775, 74401, 706, 22578, 609, 77930))
$ws.Range("B10:Q10").Value2 = (Row2D @(1649, 105558, 1594, 10241, 1434, 9183, 753, 96415, 719, 45040, 714, 137123, 601, 38448, 521, 13368))
$ws.Range("B11:Q11").Value2 = (Row2D @(1786, 228612, 961, 123053, 1351, 172973, 800, 204702, 667, 85419, 624, 239641, 657, 84050, 527, 269730))

$ws.Range("A12").Value2 = 262144
$ws.Range("B12:Q12").Value2 = (Row2D @(1761, 450942, 764, 195689, 1236, 316479, 755, 386678, 887, 226991, 338, 259424, "Failure", "Failure", "Failure", "Failure"))

$ws.Range("A13").Value2 = 524288
$ws.Range("B13:Q13").Value2 = (Row2D @(1292, 661328, 556, 284603, "Failed", "failed", "failed", "failed", "failure", "failure", "FAILURE", "failure", "failure", "failure", "failure", "failure"))

$ws.Range("A14").Formula = "=A13*2"
$ws.Range("B14:Q14").Value2 = (Row2D @("Failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed"))

# --- Second run (rows 29-41), transport/section labels then data ---
$ws.Range("A29").Value2 = "ipc:"
$ws.Range("C29").Value2 = "1 repond"
$ws.Range("G29").Value2 = "2 repond"
$ws.Range("O29").Value2 = "4 respond"
$ws.Range("K29").Value2 = "3 respond"

$ws.Range("A30:Q30").Value2 = (Row2D @("Size", "big surveys/sec", "big KBs", "small surveys/sec", "small KBs", "big surveys/sec", "big KBs", "small surveys/sec", "small KBs", "big surveys/sec", "big KBs", "small surveys/sec", "small KBs", "big surveys/sec", "big KBs", "small surveys/sec", "small KBs"))
$ws.Range("S30").Value2 = " "

$ws.Range("A31:A38").Value2 = (Col2D @(1024, 2048, 4096, 8192, 16384, 32768, 65536, 131072))

$ws.Range("B31:Q31").Value2 = (Row2D @(2223.7, 2223.7, 1840, 1840, 1441, 1441, 1223, 2445, 1012, 1022, 1027, 3082, 713, 713, 769, 2717))
$ws.Range("B32:Q32").Value2 = (Row2D @(2358, 4716, 2085, 4169, 1475, 2949, 1463, 5852, 1023, 2046, 911, 5464, 724, 1448, 736, 5336))
$ws.Range("B33:Q33").Value2 = (Row2D @(2096, 8386, 1949, 7796, 1352, 5406, 1436, 11489, 987, 3947, 799, 9594, 736, 2944, 744, 11906))
$ws.Range("B34:Q34").Value2 = (Row2D @(2039, 16310, 1767, 14134, 1423, 11386, 1425, 22805, 979, 7832, 657, 1551, 694, 5554, 688, 22007))
$ws.Range("B35:Q35").Value2 = (Row2D @(2195, 35126, 1945, 31122, 1462, 23392, 1413, 45230, 869, 13909, 801, 38465, 759, 12143, 713, 45633))
$ws.Range("B36:Q36").Value2 = (Row2D @(2400, 76794, 1612, 51596, 1380, 44174, 1350, 86416, 972, 31122, 899, 86269, 735, 23521, 632, 80931))
$ws.Range("B37:Q37").Value2 = (Row2D @(2172, 138979, 1733, 110899, 1436, 91901, 1313, 168023, 964, 61675, 816, 156581, 676, 43284, 639, 163526))
$ws.Range("B38:Q38").Value2 = (Row2D @(2012, 257442, 1076, 137708, 1323, 169335, 994, 254574, 812, 103905, 765, 259092, 648, 82907, 538, 275639))

$ws.Range("A39").Formula = "=A38*2"
$ws.Range("B39:Q39").Value2 = (Row2D @(1598, 409142, 847, 216747, 1072, 27443, 797, 407936, 663, 169750, 303, 232911, "FAILED", "FAILED", "FAILED", "FAILED"))

$ws.Range("A40").Formula = "=A39*2"
$ws.Range("B40:Q40").Value2 = (Row2D @(1147, 587425, 545, 279278, "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED", "FAILED"))

$ws.Range("A41").Formula = "=A40*2"
$ws.Range("B41:Q41").Value2 = (Row2D @("failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed", "failed"))

# Column widths for the newly used columns (J, L, N, O, P, R) - autofit to content
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666
$ws.Columns.Item(12).ColumnWidth = 16
$ws.Columns.Item(14).ColumnWidth = 16
$ws.Columns.Item(15).ColumnWidth = 11.166666666666666
$ws.Columns.Item(16).ColumnWidth = 16
$ws.Columns.Item(18).ColumnWidth = 16

# Scroll / selection state to match the saved view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q41").Select()
